$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.022883666666667
$ws.Range("H2").Value = 3.068651
$ws.Range("I2").Value = 0.1594660351460709
$ws.Range("J2").Value = 0.1594660351460709
$ws.Range("M2").Value = 15.75563966666667
$ws.Range("N2").Value = 47.266919
$ws.Range("O2").Value = 0.3220556913988901
$ws.Range("P2").Value = 0.32205569139889
$ws.Range("Q2").Value = 16.11618647291878
$ws.Range("R2").Value = 145.045678256269
$ws.Range("S2").Value = 0.05135694420360758
$ws.Range("T2").Value = 0.05135694420360756

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.022883666666667
$ws.Range("H3").Value = 3.068651
$ws.Range("I3").Value = 0.1594660351460709
$ws.Range("J3").Value = 0.1594660351460709
$ws.Range("O3").Value = 0.5509544596378365
$ws.Range("P3").Value = 0.5509544596378364
$ws.Range("Q3").Value = 27.57065019109356
$ws.Range("R3").Value = 248.135851719842
$ws.Range("S3").Value = 0.08785852322449177
$ws.Range("T3").Value = 0.08785852322449172

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.022883666666667
$ws.Range("H4").Value = 3.068651
$ws.Range("I4").Value = 0.1594660351460709
$ws.Range("J4").Value = 0.1594660351460709
$ws.Range("O4").Value = 0.1269898489632735
$ws.Range("P4").Value = 0.1269898489632735
$ws.Range("Q4").Value = 6.354776955408778
$ws.Range("R4").Value = 57.192992598679
$ws.Range("S4").Value = 0.02025056771797162
$ws.Range("T4").Value = 0.02025056771797161

# Row 5
$ws.Range("I5").Value = 0.1603506552336246
$ws.Range("J5").Value = 0.1603506552336246
$ws.Range("M5").Value = 15.75563966666667
$ws.Range("N5").Value = 47.266919
$ws.Range("O5").Value = 0.3220556913988901
$ws.Range("P5").Value = 0.32205569139889
$ws.Range("Q5").Value = 16.20558922426734
$ws.Range("R5").Value = 145.850303018406
$ws.Range("S5").Value = 0.05164184113753001
$ws.Range("T5").Value = 0.05164184113753

# Row 6
$ws.Range("I6").Value = 0.1603506552336246
$ws.Range("J6").Value = 0.1603506552336246
$ws.Range("O6").Value = 0.5509544596378365
$ws.Range("P6").Value = 0.5509544596378364
$ws.Range("S6").Value = 0.08834590860681465
$ws.Range("T6").Value = 0.08834590860681463

# Row 7
$ws.Range("I7").Value = 0.1603506552336246
$ws.Range("J7").Value = 0.1603506552336246
$ws.Range("O7").Value = 0.1269898489632735
$ws.Range("P7").Value = 0.1269898489632735
$ws.Range("S7").Value = 0.02036290548927993
$ws.Range("T7").Value = 0.02036290548927993

# Row 8
$ws.Range("I8").Value = 0.6801833096203046
$ws.Range("J8").Value = 0.6801833096203045
$ws.Range("M8").Value = 15.75563966666667
$ws.Range("N8").Value = 47.266919
$ws.Range("O8").Value = 0.3220556913988901
$ws.Range("P8").Value = 0.32205569139889
$ws.Range("Q8").Value = 68.74166679799067
$ws.Range("R8").Value = 618.6750011819161
$ws.Range("S8").Value = 0.2190569060577525
$ws.Range("T8").Value = 0.2190569060577524

# Row 9
$ws.Range("I9").Value = 0.6801833096203046
$ws.Range("J9").Value = 0.6801833096203045
$ws.Range("O9").Value = 0.5509544596378365
$ws.Range("P9").Value = 0.5509544596378364
$ws.Range("S9").Value = 0.3747500278065302
$ws.Range("T9").Value = 0.3747500278065301

# Row 10
$ws.Range("I10").Value = 0.6801833096203046
$ws.Range("J10").Value = 0.6801833096203045
$ws.Range("O10").Value = 0.1269898489632735
$ws.Range("P10").Value = 0.1269898489632735
$ws.Range("S10").Value = 0.086376375756022
$ws.Range("T10").Value = 0.08637637575602196
